# Indexant llibre matrimonis 1749-1770 fins pag 30
# Adds 54 new marriage-index rows (714-767) to Hoja1, matching the
# historical write order so the shared-string table lines up exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: columns filled in bulk across the whole new block first ---
# Column M ("Anys") - constant value for this whole batch of entries
$ws.Range("M714:M767").Value = "1749-1770"

# Column L ("Serie") - constant value for this whole batch of entries
$ws.Range("L714:L767").Value = "A,2"

# --- Step 2: column C ("Cognoms Familia") filled row by row; row 719 was
# corrected after rows 720-722 were already entered, so it is written last
# among that cluster to reproduce the original shared-string allocation ---
$ws.Cells.Item(714, 3).Value = "Mas Martí"
$ws.Cells.Item(715, 3).Value = "Incognit Planes"
$ws.Cells.Item(716, 3).Value = "Civit Incognit"
$ws.Cells.Item(717, 3).Value = "Marti Ginesta"
$ws.Cells.Item(718, 3).Value = "Pedrós Rey"
$ws.Cells.Item(720, 3).Value = "Puig Valles"
$ws.Cells.Item(721, 3).Value = "Martí Solé"
$ws.Cells.Item(722, 3).Value = "Corberó Thomas"
$ws.Cells.Item(719, 3).Value = "Mas Planes"
$ws.Cells.Item(723, 3).Value = "Tarragó Solà"
$ws.Cells.Item(724, 3).Value = "Gene Arrufat"
$ws.Cells.Item(725, 3).Value = "Martí Mas"
$ws.Cells.Item(726, 3).Value = "Batlle Vergè"
$ws.Cells.Item(727, 3).Value = "Civit Rossell"
$ws.Cells.Item(728, 3).Value = "Senallosa Domingo"
$ws.Cells.Item(729, 3).Value = "Vilaplana Vergé"
$ws.Cells.Item(730, 3).Value = "Solà Mas"
$ws.Cells.Item(731, 3).Value = "Cascalló Corberó"
$ws.Cells.Item(732, 3).Value = "Duch Sucarrat"
$ws.Cells.Item(733, 3).Value = "Gene Palou"
$ws.Cells.Item(734, 3).Value = "Pedrós Felip"
$ws.Cells.Item(735, 3).Value = "Trepat Valles"
$ws.Cells.Item(736, 3).Value = "Martí Roma"
$ws.Cells.Item(737, 3).Value = "Agulló Cases"
$ws.Cells.Item(738, 3).Value = "Font Fortuny"
$ws.Cells.Item(739, 3).Value = "Cascalló Ribes"
$ws.Cells.Item(740, 3).Value = "Queral Vergé"
$ws.Cells.Item(741, 3).Value = "Mata Senallosa"
$ws.Cells.Item(742, 3).Value = "Cuberes Solsona"
$ws.Cells.Item(743, 3).Value = "Tarroja Pujades"
$ws.Cells.Item(744, 3).Value = "Falcó Mas"
$ws.Cells.Item(745, 3).Value = "Falcó Torruella"
$ws.Cells.Item(746, 3).Value = "Vergé Galceran"
$ws.Cells.Item(747, 3).Value = "Pedrós Mas"
$ws.Cells.Item(748, 3).Value = "Pallas Pujal"
$ws.Cells.Item(749, 3).Value = "Mas Torren"
$ws.Cells.Item(750, 3).Value = "Dalmau Vallés"
$ws.Cells.Item(751, 3).Value = "Coll Fabregat"
$ws.Cells.Item(752, 3).Value = "Mosset Rubiol"
$ws.Cells.Item(753, 3).Value = "Vergé Tapies"
$ws.Cells.Item(754, 3).Value = "Civit Serra"
$ws.Cells.Item(755, 3).Value = "Vergé Mas"
$ws.Cells.Item(756, 3).Value = "Curcó Pujol"
$ws.Cells.Item(757, 3).Value = "Pujol Roigé"
$ws.Cells.Item(758, 3).Value = "Cascalló Druet"
$ws.Cells.Item(759, 3).Value = "Mas Roige"
$ws.Cells.Item(760, 3).Value = "Pujol Prats"
$ws.Cells.Item(761, 3).Value = "Salvat Bonjorn"
$ws.Cells.Item(762, 3).Value = "Torra Balagué"
$ws.Cells.Item(763, 3).Value = "Sucarrat Pera"
$ws.Cells.Item(764, 3).Value = "Mosset Codina"
$ws.Cells.Item(765, 3).Value = "Ponsarnau Palou"
$ws.Cells.Item(766, 3).Value = "Batalla Tarroja"
$ws.Cells.Item(767, 3).Value = "vergé Mas"

# --- Step 3: remaining columns (A, B, J, K) - no new shared strings ---
$ws.Cells.Item(714, 1).Value = 1749
$ws.Cells.Item(714, 2).Value = 3
$ws.Cells.Item(714, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(714, 11).Value = 47
$ws.Cells.Item(715, 1).Value = 1750
$ws.Cells.Item(715, 2).Value = 3
$ws.Cells.Item(715, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(715, 11).Value = 47
$ws.Cells.Item(716, 1).Value = 1750
$ws.Cells.Item(716, 2).Value = 4
$ws.Cells.Item(716, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(716, 11).Value = 47
$ws.Cells.Item(717, 1).Value = 1750
$ws.Cells.Item(717, 2).Value = 4
$ws.Cells.Item(717, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(717, 11).Value = 47
$ws.Cells.Item(718, 1).Value = 1751
$ws.Cells.Item(718, 2).Value = 5
$ws.Cells.Item(718, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(718, 11).Value = 47
$ws.Cells.Item(719, 1).Value = 1751
$ws.Cells.Item(719, 2).Value = 5
$ws.Cells.Item(719, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(719, 11).Value = 47
$ws.Cells.Item(720, 1).Value = 1751
$ws.Cells.Item(720, 2).Value = 5
$ws.Cells.Item(720, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(720, 11).Value = 47
$ws.Cells.Item(721, 1).Value = 1753
$ws.Cells.Item(721, 2).Value = 6
$ws.Cells.Item(721, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(721, 11).Value = 47
$ws.Cells.Item(722, 1).Value = 1753
$ws.Cells.Item(722, 2).Value = 6
$ws.Cells.Item(722, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(722, 11).Value = 47
$ws.Cells.Item(723, 1).Value = 1753
$ws.Cells.Item(723, 2).Value = 7
$ws.Cells.Item(723, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(723, 11).Value = 47
$ws.Cells.Item(724, 1).Value = 1753
$ws.Cells.Item(724, 2).Value = 7
$ws.Cells.Item(724, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(724, 11).Value = 47
$ws.Cells.Item(725, 1).Value = 1754
$ws.Cells.Item(725, 2).Value = 7
$ws.Cells.Item(725, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(725, 11).Value = 47
$ws.Cells.Item(726, 1).Value = 1754
$ws.Cells.Item(726, 2).Value = 8
$ws.Cells.Item(726, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(726, 11).Value = 47
$ws.Cells.Item(727, 1).Value = 1755
$ws.Cells.Item(727, 2).Value = 9
$ws.Cells.Item(727, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(727, 11).Value = 47
$ws.Cells.Item(728, 1).Value = 1755
$ws.Cells.Item(728, 2).Value = 9
$ws.Cells.Item(728, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(728, 11).Value = 47
$ws.Cells.Item(729, 1).Value = 1755
$ws.Cells.Item(729, 2).Value = 9
$ws.Cells.Item(729, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(729, 11).Value = 47
$ws.Cells.Item(730, 1).Value = 1755
$ws.Cells.Item(730, 2).Value = 10
$ws.Cells.Item(730, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(730, 11).Value = 47
$ws.Cells.Item(731, 1).Value = 1755
$ws.Cells.Item(731, 2).Value = 10
$ws.Cells.Item(731, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(731, 11).Value = 47
$ws.Cells.Item(732, 1).Value = 1756
$ws.Cells.Item(732, 2).Value = 10
$ws.Cells.Item(732, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(732, 11).Value = 47
$ws.Cells.Item(733, 1).Value = 1756
$ws.Cells.Item(733, 2).Value = 11
$ws.Cells.Item(733, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(733, 11).Value = 47
$ws.Cells.Item(734, 1).Value = 1756
$ws.Cells.Item(734, 2).Value = 11
$ws.Cells.Item(734, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(734, 11).Value = 47
$ws.Cells.Item(735, 1).Value = 1756
$ws.Cells.Item(735, 2).Value = 12
$ws.Cells.Item(735, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(735, 11).Value = 47
$ws.Cells.Item(736, 1).Value = 1756
$ws.Cells.Item(736, 2).Value = 12
$ws.Cells.Item(736, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(736, 11).Value = 47
$ws.Cells.Item(737, 1).Value = 1756
$ws.Cells.Item(737, 2).Value = 13
$ws.Cells.Item(737, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(737, 11).Value = 47
$ws.Cells.Item(738, 1).Value = 1756
$ws.Cells.Item(738, 2).Value = 13
$ws.Cells.Item(738, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(738, 11).Value = 47
$ws.Cells.Item(739, 1).Value = 1756
$ws.Cells.Item(739, 2).Value = 14
$ws.Cells.Item(739, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(739, 11).Value = 47
$ws.Cells.Item(740, 1).Value = 1757
$ws.Cells.Item(740, 2).Value = 14
$ws.Cells.Item(740, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(740, 11).Value = 47
$ws.Cells.Item(741, 1).Value = 1757
$ws.Cells.Item(741, 2).Value = 15
$ws.Cells.Item(741, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(741, 11).Value = 47
$ws.Cells.Item(742, 1).Value = 1757
$ws.Cells.Item(742, 2).Value = 15
$ws.Cells.Item(742, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(742, 11).Value = 47
$ws.Cells.Item(743, 1).Value = 1757
$ws.Cells.Item(743, 2).Value = 16
$ws.Cells.Item(743, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(743, 11).Value = 47
$ws.Cells.Item(744, 1).Value = 1757
$ws.Cells.Item(744, 2).Value = 16
$ws.Cells.Item(744, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(744, 11).Value = 47
$ws.Cells.Item(745, 1).Value = 1758
$ws.Cells.Item(745, 2).Value = 17
$ws.Cells.Item(745, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(745, 11).Value = 47
$ws.Cells.Item(746, 1).Value = 1758
$ws.Cells.Item(746, 2).Value = 17
$ws.Cells.Item(746, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(746, 11).Value = 47
$ws.Cells.Item(747, 1).Value = 1758
$ws.Cells.Item(747, 2).Value = 18
$ws.Cells.Item(747, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(747, 11).Value = 47
$ws.Cells.Item(748, 1).Value = 1758
$ws.Cells.Item(748, 2).Value = 18
$ws.Cells.Item(748, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(748, 11).Value = 47
$ws.Cells.Item(749, 1).Value = 1758
$ws.Cells.Item(749, 2).Value = 19
$ws.Cells.Item(749, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(749, 11).Value = 47
$ws.Cells.Item(750, 1).Value = 1758
$ws.Cells.Item(750, 2).Value = 19
$ws.Cells.Item(750, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(750, 11).Value = 47
$ws.Cells.Item(751, 1).Value = 1758
$ws.Cells.Item(751, 2).Value = 20
$ws.Cells.Item(751, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(751, 11).Value = 47
$ws.Cells.Item(752, 1).Value = 1758
$ws.Cells.Item(752, 2).Value = 20
$ws.Cells.Item(752, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(752, 11).Value = 47
$ws.Cells.Item(753, 1).Value = 1758
$ws.Cells.Item(753, 2).Value = 21
$ws.Cells.Item(753, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(753, 11).Value = 47
$ws.Cells.Item(754, 1).Value = 1759
$ws.Cells.Item(754, 2).Value = 22
$ws.Cells.Item(754, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(754, 11).Value = 47
$ws.Cells.Item(755, 1).Value = 1759
$ws.Cells.Item(755, 2).Value = 22
$ws.Cells.Item(755, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(755, 11).Value = 47
$ws.Cells.Item(756, 1).Value = 1759
$ws.Cells.Item(756, 2).Value = 23
$ws.Cells.Item(756, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(756, 11).Value = 47
$ws.Cells.Item(757, 1).Value = 1759
$ws.Cells.Item(757, 2).Value = 23
$ws.Cells.Item(757, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(757, 11).Value = 47
$ws.Cells.Item(758, 1).Value = 1759
$ws.Cells.Item(758, 2).Value = 24
$ws.Cells.Item(758, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(758, 11).Value = 47
$ws.Cells.Item(759, 1).Value = 1759
$ws.Cells.Item(759, 2).Value = 24
$ws.Cells.Item(759, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(759, 11).Value = 47
$ws.Cells.Item(760, 1).Value = 1759
$ws.Cells.Item(760, 2).Value = 25
$ws.Cells.Item(760, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(760, 11).Value = 47
$ws.Cells.Item(761, 1).Value = 1759
$ws.Cells.Item(761, 2).Value = 25
$ws.Cells.Item(761, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(761, 11).Value = 47
$ws.Cells.Item(762, 1).Value = 1760
$ws.Cells.Item(762, 2).Value = 26
$ws.Cells.Item(762, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(762, 11).Value = 47
$ws.Cells.Item(763, 1).Value = 1761
$ws.Cells.Item(763, 2).Value = 28
$ws.Cells.Item(763, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(763, 11).Value = 47
$ws.Cells.Item(764, 1).Value = 1761
$ws.Cells.Item(764, 2).Value = 28
$ws.Cells.Item(764, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(764, 11).Value = 47
$ws.Cells.Item(765, 1).Value = 1761
$ws.Cells.Item(765, 2).Value = 29
$ws.Cells.Item(765, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(765, 11).Value = 47
$ws.Cells.Item(766, 1).Value = 1761
$ws.Cells.Item(766, 2).Value = 29
$ws.Cells.Item(766, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(766, 11).Value = 47
$ws.Cells.Item(767, 1).Value = 1761
$ws.Cells.Item(767, 2).Value = 30
$ws.Cells.Item(767, 10).Value = "SPN 2,02 C"
$ws.Cells.Item(767, 11).Value = 47

# --- Step 4: extend the AutoFilter range and the _FilterDatabase name ---
$ws.AutoFilterMode = $false
$ws.Range("A1:M715").AutoFilter()
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "=Hoja1!`$A`$1:`$M`$715"

# --- Step 5: move the active selection to the next empty row, as left by
# the original author after finishing data entry ---
$ws.Range("A768").Select()

